# FBOW_Sensitivity 2024 Edition update
# - renames the FLOW_Sensitivity sheet to FBOW_Sensitivity
# - updates the Fixed Charge Rate input assumptions (rows 7-9) on the
#   "Fixed Charge Rates" sheet
# - re-points the F5 "base value" cell on FBOW_Sensitivity at the newly
#   recomputed FCR (nominal) coming off the Fixed Charge Rates sheet
# - applies light number formatting to a few of the key input cells
# All downstream formulas (B2:D6 on FBOW_Sensitivity, rows 10-19 on
# Fixed Charge Rates, etc.) are plain formulas already in the workbook, so
# they recalc automatically once the inputs below change.

$wb = $excel.ActiveWorkbook

# --- rename the sheet (formulas on "Fixed Charge Rates" that reference it
#     follow automatically) ---------------------------------------------
$wsFbow = $wb.Worksheets.Item("FLOW_Sensitivity")
$wsFbow.Name = "FBOW_Sensitivity"

$wsFcr = $wb.Worksheets.Item("Fixed Charge Rates")

# --- updated Fixed Charge Rate assumptions (row 7: Debt fraction,
#     row 8: Debt interest rate, row 9: Return on equity) ----------------
$wsFcr.Range("C7:G7").Value = 0.734
$wsFcr.Range("C8:G8").Value = 0.07
$wsFcr.Range("C9:G9").Value = 0.105

# --- FBOW_Sensitivity base CapEx/OpEx/etc cells keep their values, just
#     pick up light number formatting ------------------------------------
$wsFbow.Range("F2").NumberFormat = "0"
$wsFbow.Range("F3").NumberFormat = "0"
$wsFbow.Range("F4").NumberFormat = "0.0"
$wsFbow.Range("F5").NumberFormat = "0.00"

# --- base WACC (nominal) cell now pulls straight from the recomputed FCR
#     (nominal) on the Fixed Charge Rates sheet instead of being a typed
#     constant ------------------------------------------------------------
$wsFbow.Range("F5").Formula = "='Fixed Charge Rates'!C10*100"

# --- restore the on-screen selections recorded in the saved workbook, and
#     make sure FBOW_Sensitivity is the active (left-most/visible) tab ----
$wsFcr.Range("H6").Select()
$wsFbow.Range("F6").Select()
$wsFbow.Activate()
